$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new survey rows (94-102) covering survey waves 31-35 ---
$ws.Cells.Item(94,1).Value = 3
$ws.Cells.Item(94,2).Value = 0
$ws.Cells.Item(94,3).Value = "uk"
$ws.Cells.Item(94,4).Value = 81
$ws.Cells.Item(94,5).Value = "F"
$ws.Cells.Item(94,6).Value = 31
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(94,7))
$ws.Cells.Item(94,7).Value = 44483
$ws.Cells.Item(94,8).Value = "21-037554_PFW31_Final_ICUO"

$ws.Cells.Item(95,1).Value = 3
$ws.Cells.Item(95,3).Value = "uk"
$ws.Cells.Item(95,4).Value = 82
$ws.Cells.Item(95,5).Value = "E"
$ws.Cells.Item(95,6).Value = 32
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(95,7))
$ws.Cells.Item(95,7).Value = 44491
$ws.Cells.Item(95,8).Value = "21-037558_PEW32_Final_ICUO"

$ws.Cells.Item(96,1).Value = 3
$ws.Cells.Item(96,3).Value = "uk"
$ws.Cells.Item(96,4).Value = 83
$ws.Cells.Item(96,5).Value = "F"
$ws.Cells.Item(96,6).Value = 32
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(96,7))
$ws.Cells.Item(96,7).Value = 44498
$ws.Cells.Item(96,8).Value = "21-037554_PFW32_Final_ICUO"

$ws.Cells.Item(97,1).Value = 3
$ws.Cells.Item(97,3).Value = "uk"
$ws.Cells.Item(97,4).Value = 84
$ws.Cells.Item(97,5).Value = "E"
$ws.Cells.Item(97,6).Value = 33
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(97,7))
$ws.Cells.Item(97,7).Value = 44504
$ws.Cells.Item(97,8).Value = "21-088043_PEW33_Final_ICUO"

$ws.Cells.Item(98,1).Value = 3
$ws.Cells.Item(98,3).Value = "uk"
$ws.Cells.Item(98,4).Value = 85
$ws.Cells.Item(98,5).Value = "F"
$ws.Cells.Item(98,6).Value = 33
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(98,7))
$ws.Cells.Item(98,7).Value = 44512
$ws.Cells.Item(98,8).Value = "21-088071_PFW33_Final_ICUO"

$ws.Cells.Item(99,1).Value = 3
$ws.Cells.Item(99,3).Value = "uk"
$ws.Cells.Item(99,4).Value = 86
$ws.Cells.Item(99,5).Value = "E"
$ws.Cells.Item(99,6).Value = 34
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(99,7))
$ws.Cells.Item(99,7).Value = 44519
$ws.Cells.Item(99,8).Value = "21-088043_PEW34_Final_ICUO"

$ws.Cells.Item(100,1).Value = 3
$ws.Cells.Item(100,3).Value = "uk"
$ws.Cells.Item(100,4).Value = 87
$ws.Cells.Item(100,5).Value = "F"
$ws.Cells.Item(100,6).Value = 34
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(100,7))
$ws.Cells.Item(100,7).Value = 44525
$ws.Cells.Item(100,8).Value = "21-088071_PFW34_Final_AMENDED_IntUse"

$ws.Cells.Item(101,1).Value = 3
$ws.Cells.Item(101,3).Value = "uk"
$ws.Cells.Item(101,4).Value = 88
$ws.Cells.Item(101,5).Value = "E"
$ws.Cells.Item(101,6).Value = 35
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(101,7))
$ws.Cells.Item(101,7).Value = 44533
$ws.Cells.Item(101,8).Value = "21-088043_PEW35_Final_ICUO"

$ws.Cells.Item(102,1).Value = 3
$ws.Cells.Item(102,3).Value = "uk"
$ws.Cells.Item(102,4).Value = 89
$ws.Cells.Item(102,5).Value = "F"
$ws.Cells.Item(102,6).Value = 35
$ws.Cells.Item(2,7).Copy($ws.Cells.Item(102,7))
$ws.Cells.Item(102,7).Value = 44540
$ws.Cells.Item(102,8).Value = "21-088071_PFW35_Final_ICUO"

# Recreate the "spss_name" lookup formula (shared across I83:I102) for the new rows
$ws.Range("I94:I102").Formula = '=C94&"_"&"sr"&TEXT(D94,"00")&"_"&YEAR(G94)&TEXT(G94,"MM")&TEXT(G94,"DD")&"_p"&E94&"_wv"&TEXT(F94,"00")&""'

# Move the saved selection/scroll position to reflect where editing left off
$ws.Range("H116").Select()
